# ------------------------------------------------------------------
# Edit: raven.docx
#   1. Append two trailing spaces to the first paragraph's existing
#      text, then append three additional red-colored runs that spell
#      out "(This is a change – Version for main branch)" (split
#      across three runs exactly as the source OOXML does).
#   2. Delete the trailing paragraph that reads
#      "ank God almighty, we are free at last."
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. First paragraph: extend + add red commentary runs ----------
$p1 = $d.Paragraphs.Item(1).Range

# Cursor = content position immediately before the paragraph mark.
$cursor = $p1.End - 1

# Run 1 (no new formatting, just two trailing spaces appended to the
# existing black run's text).
$run1Text = "  "
$p1.InsertAfter($run1Text)
$cursor = $cursor + $run1Text.Length

# Run 2 (red)
$run2Text = "(This is a change " + [char]0x2013 + " Ve"
$run2Start = $cursor
$p1.InsertAfter($run2Text)
$cursor = $cursor + $run2Text.Length
$run2End = $cursor
$d.Range($run2Start, $run2End).Font.Color = 255

# Run 3 (red)
$run3Text = "rsion for main branch"
$run3Start = $cursor
$p1.InsertAfter($run3Text)
$cursor = $cursor + $run3Text.Length
$run3End = $cursor
$d.Range($run3Start, $run3End).Font.Color = 255

# Run 4 (red)
$run4Text = ")"
$run4Start = $cursor
$p1.InsertAfter($run4Text)
$cursor = $cursor + $run4Text.Length
$run4End = $cursor
$d.Range($run4Start, $run4End).Font.Color = 255

# --- 2. Remove the trailing "...we are free at last." paragraph ----
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastRange = $lastPara.Range
$lastRange.Delete()
